$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44740
$ws.Cells.Item(2, 11).Value = 2500
$ws.Cells.Item(2, 12).Value = 2500
$ws.Cells.Item(2, 13).Value = 2500
$ws.Cells.Item(2, 16).Value = 2500

# Row 3
$ws.Cells.Item(3, 4).Value = 44749
$ws.Cells.Item(3, 10).Value = 80

# Row 4
$ws.Cells.Item(4, 4).Value = 44453
$ws.Cells.Item(4, 10).Value = 20
$ws.Cells.Item(4, 11).Value = 2300
$ws.Cells.Item(4, 12).Value = 2300
$ws.Cells.Item(4, 13).Value = 2300
$ws.Cells.Item(4, 16).Value = 2300

# Row 5
$ws.Cells.Item(5, 4).Value = 44910
$ws.Cells.Item(5, 11).Value = 2500
$ws.Cells.Item(5, 12).Value = 2500
$ws.Cells.Item(5, 13).Value = 2500
$ws.Cells.Item(5, 16).Value = 2500

# Row 6
$ws.Cells.Item(6, 4).Value = 44837
$ws.Cells.Item(6, 10).Value = 50

# Row 7
$ws.Cells.Item(7, 4).Value = 44677
$ws.Cells.Item(7, 11).Value = 5500
$ws.Cells.Item(7, 12).Value = 5500
$ws.Cells.Item(7, 13).Value = 5500
$ws.Cells.Item(7, 16).Value = 5500

# Row 8
$ws.Cells.Item(8, 4).Value = 44895
$ws.Cells.Item(8, 10).Value = 40

# Row 9
$ws.Cells.Item(9, 4).Value = 44781
$ws.Cells.Item(9, 10).Value = 250
$ws.Cells.Item(9, 11).Value = 2700
$ws.Cells.Item(9, 12).Value = 2700
$ws.Cells.Item(9, 13).Value = 2700
$ws.Cells.Item(9, 16).Value = 2700

# Row 10
$ws.Cells.Item(10, 4).Value = 44474
$ws.Cells.Item(10, 10).Value = 20
$ws.Cells.Item(10, 11).Value = 1600
$ws.Cells.Item(10, 12).Value = 1600
$ws.Cells.Item(10, 13).Value = 1600
$ws.Cells.Item(10, 16).Value = 1600

# Row 11
$ws.Cells.Item(11, 4).Value = 44930
$ws.Cells.Item(11, 10).Value = 90

# Row 12
$ws.Cells.Item(12, 4).Value = 44769
$ws.Cells.Item(12, 10).Value = 140
$ws.Cells.Item(12, 11).Value = 3300
$ws.Cells.Item(12, 12).Value = 3300
$ws.Cells.Item(12, 13).Value = 3300
$ws.Cells.Item(12, 16).Value = 3300

# Row 13
$ws.Cells.Item(13, 4).Value = 44783
$ws.Cells.Item(13, 10).Value = 90
$ws.Cells.Item(13, 11).Value = 2700
$ws.Cells.Item(13, 12).Value = 2700
$ws.Cells.Item(13, 13).Value = 2700
$ws.Cells.Item(13, 16).Value = 2700

# Row 14
$ws.Cells.Item(14, 4).Value = 44771
$ws.Cells.Item(14, 10).Value = 30
$ws.Cells.Item(14, 11).Value = 3300
$ws.Cells.Item(14, 12).Value = 3300
$ws.Cells.Item(14, 13).Value = 3300
$ws.Cells.Item(14, 16).Value = 3300

# Row 15
$ws.Cells.Item(15, 4).Value = 44811
$ws.Cells.Item(15, 10).Value = 120
$ws.Cells.Item(15, 11).Value = 2700
$ws.Cells.Item(15, 12).Value = 2700
$ws.Cells.Item(15, 13).Value = 2700
$ws.Cells.Item(15, 16).Value = 2700

# Row 16
$ws.Cells.Item(16, 4).Value = 44839
$ws.Cells.Item(16, 10).Value = 80
$ws.Cells.Item(16, 11).Value = 2500
$ws.Cells.Item(16, 12).Value = 2500
$ws.Cells.Item(16, 13).Value = 2500
$ws.Cells.Item(16, 16).Value = 2500

# Row 17
$ws.Cells.Item(17, 4).Value = 44685
$ws.Cells.Item(17, 10).Value = 60
$ws.Cells.Item(17, 11).Value = 5000
$ws.Cells.Item(17, 12).Value = 6000
$ws.Cells.Item(17, 13).Value = 5333
$ws.Cells.Item(17, 16).Value = 5333

# Row 18
$ws.Cells.Item(18, 4).Value = 44797
$ws.Cells.Item(18, 10).Value = 200
$ws.Cells.Item(18, 11).Value = 2700
$ws.Cells.Item(18, 12).Value = 2700
$ws.Cells.Item(18, 13).Value = 2700
$ws.Cells.Item(18, 16).Value = 2700

# Row 19
$ws.Cells.Item(19, 4).Value = 44706
$ws.Cells.Item(19, 10).Value = 90
$ws.Cells.Item(19, 11).Value = 4700
$ws.Cells.Item(19, 12).Value = 4700
$ws.Cells.Item(19, 13).Value = 4700
$ws.Cells.Item(19, 16).Value = 4700

# Row 20
$ws.Cells.Item(20, 4).Value = 44893
$ws.Cells.Item(20, 10).Value = 50
$ws.Cells.Item(20, 11).Value = 2500
$ws.Cells.Item(20, 12).Value = 2500
$ws.Cells.Item(20, 13).Value = 2500
$ws.Cells.Item(20, 16).Value = 2500

# Row 21
$ws.Cells.Item(21, 4).Value = 44452
$ws.Cells.Item(21, 10).Value = 120
$ws.Cells.Item(21, 11).Value = 2300
$ws.Cells.Item(21, 12).Value = 2300
$ws.Cells.Item(21, 13).Value = 2300
$ws.Cells.Item(21, 16).Value = 2300

# Row 22
$ws.Cells.Item(22, 4).Value = 44669
$ws.Cells.Item(22, 10).Value = 60
$ws.Cells.Item(22, 11).Value = 6250
$ws.Cells.Item(22, 12).Value = 6250
$ws.Cells.Item(22, 13).Value = 6250
$ws.Cells.Item(22, 16).Value = 6250

# Row 23
$ws.Cells.Item(23, 4).Value = 44767
$ws.Cells.Item(23, 10).Value = 180
$ws.Cells.Item(23, 11).Value = 3300
$ws.Cells.Item(23, 12).Value = 3300
$ws.Cells.Item(23, 13).Value = 3300
$ws.Cells.Item(23, 16).Value = 3300

# Row 24
$ws.Cells.Item(24, 4).Value = 44497
$ws.Cells.Item(24, 10).Value = 50

# Row 25
$ws.Cells.Item(25, 4).Value = 44816
$ws.Cells.Item(25, 11).Value = 2700
$ws.Cells.Item(25, 12).Value = 2700
$ws.Cells.Item(25, 13).Value = 2700
$ws.Cells.Item(25, 16).Value = 2700

# Row 26
$ws.Cells.Item(26, 4).Value = 44720
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 3600
$ws.Cells.Item(26, 12).Value = 3600
$ws.Cells.Item(26, 13).Value = 3600
$ws.Cells.Item(26, 16).Value = 3600

# Row 27
$ws.Cells.Item(27, 4).Value = 44473
$ws.Cells.Item(27, 10).Value = 140
$ws.Cells.Item(27, 11).Value = 1600
$ws.Cells.Item(27, 12).Value = 1600
$ws.Cells.Item(27, 13).Value = 1600
$ws.Cells.Item(27, 16).Value = 1600

# Row 28
$ws.Cells.Item(28, 4).Value = 44496
$ws.Cells.Item(28, 10).Value = 40

# Row 29
$ws.Cells.Item(29, 4).Value = 44679
$ws.Cells.Item(29, 10).Value = 30
$ws.Cells.Item(29, 11).Value = 5500
$ws.Cells.Item(29, 12).Value = 5500
$ws.Cells.Item(29, 13).Value = 5500
$ws.Cells.Item(29, 16).Value = 5500

# Row 30
$ws.Cells.Item(30, 4).Value = 44203
$ws.Cells.Item(30, 10).Value = 30
$ws.Cells.Item(30, 11).Value = 2000
$ws.Cells.Item(30, 12).Value = 2000
$ws.Cells.Item(30, 13).Value = 2000
$ws.Cells.Item(30, 16).Value = 2000

# Row 31
$ws.Cells.Item(31, 4).Value = 44753
$ws.Cells.Item(31, 10).Value = 130
$ws.Cells.Item(31, 12).Value = 3300
$ws.Cells.Item(31, 13).Value = 2931
$ws.Cells.Item(31, 16).Value = 2931

# Row 32
$ws.Cells.Item(32, 4).Value = 44868
$ws.Cells.Item(32, 10).Value = 80

# Row 33
$ws.Cells.Item(33, 4).Value = 44741
$ws.Cells.Item(33, 10).Value = 100
$ws.Cells.Item(33, 11).Value = 2500
$ws.Cells.Item(33, 12).Value = 2500
$ws.Cells.Item(33, 13).Value = 2500
$ws.Cells.Item(33, 16).Value = 2500

# Row 34
$ws.Cells.Item(34, 4).Value = 44755
$ws.Cells.Item(34, 10).Value = 90

# Row 35
$ws.Cells.Item(35, 4).Value = 44776
$ws.Cells.Item(35, 10).Value = 100
$ws.Cells.Item(35, 11).Value = 2700
$ws.Cells.Item(35, 12).Value = 2700
$ws.Cells.Item(35, 13).Value = 2700
$ws.Cells.Item(35, 16).Value = 2700

# Row 36
$ws.Cells.Item(36, 4).Value = 44756
$ws.Cells.Item(36, 10).Value = 120
$ws.Cells.Item(36, 11).Value = 3300
$ws.Cells.Item(36, 13).Value = 3300
$ws.Cells.Item(36, 16).Value = 3300

# Row 37
$ws.Cells.Item(37, 4).Value = 44487
$ws.Cells.Item(37, 10).Value = 50
$ws.Cells.Item(37, 11).Value = 2200
$ws.Cells.Item(37, 12).Value = 2200
$ws.Cells.Item(37, 13).Value = 2200
$ws.Cells.Item(37, 16).Value = 2200

# Row 38
$ws.Cells.Item(38, 4).Value = 44879
$ws.Cells.Item(38, 10).Value = 200
$ws.Cells.Item(38, 11).Value = 2500
$ws.Cells.Item(38, 12).Value = 2500
$ws.Cells.Item(38, 13).Value = 2500
$ws.Cells.Item(38, 16).Value = 2500

# Row 39
$ws.Cells.Item(39, 4).Value = 44747
$ws.Cells.Item(39, 10).Value = 80
$ws.Cells.Item(39, 11).Value = 2500
$ws.Cells.Item(39, 12).Value = 2500
$ws.Cells.Item(39, 13).Value = 2500
$ws.Cells.Item(39, 16).Value = 2500

# Row 40
$ws.Cells.Item(40, 4).Value = 44447
$ws.Cells.Item(40, 10).Value = 75
$ws.Cells.Item(40, 11).Value = 2200
$ws.Cells.Item(40, 12).Value = 2200
$ws.Cells.Item(40, 13).Value = 2200
$ws.Cells.Item(40, 16).Value = 2200

# Row 41
$ws.Cells.Item(41, 4).Value = 44483
$ws.Cells.Item(41, 10).Value = 50
$ws.Cells.Item(41, 11).Value = 2200
$ws.Cells.Item(41, 12).Value = 2200
$ws.Cells.Item(41, 13).Value = 2200
$ws.Cells.Item(41, 16).Value = 2200

# Row 42
$ws.Cells.Item(42, 4).Value = 44804
$ws.Cells.Item(42, 10).Value = 100
$ws.Cells.Item(42, 11).Value = 3300
$ws.Cells.Item(42, 12).Value = 3300
$ws.Cells.Item(42, 13).Value = 3300
$ws.Cells.Item(42, 16).Value = 3300

# Row 43
$ws.Cells.Item(43, 4).Value = 44809
$ws.Cells.Item(43, 10).Value = 150
$ws.Cells.Item(43, 11).Value = 2700
$ws.Cells.Item(43, 12).Value = 2700
$ws.Cells.Item(43, 13).Value = 2700
$ws.Cells.Item(43, 16).Value = 2700

# Row 44
$ws.Cells.Item(44, 4).Value = 44931
$ws.Cells.Item(44, 11).Value = 2500
$ws.Cells.Item(44, 12).Value = 2500
$ws.Cells.Item(44, 13).Value = 2500
$ws.Cells.Item(44, 16).Value = 2500

# Row 45
$ws.Cells.Item(45, 4).Value = 44719
$ws.Cells.Item(45, 11).Value = 3600
$ws.Cells.Item(45, 12).Value = 3600
$ws.Cells.Item(45, 13).Value = 3600
$ws.Cells.Item(45, 16).Value = 3600

# Row 46
$ws.Cells.Item(46, 4).Value = 44825
$ws.Cells.Item(46, 10).Value = 30
$ws.Cells.Item(46, 11).Value = 2700
$ws.Cells.Item(46, 12).Value = 2700
$ws.Cells.Item(46, 13).Value = 2700
$ws.Cells.Item(46, 16).Value = 2700

# Row 47
$ws.Cells.Item(47, 4).Value = 44484
$ws.Cells.Item(47, 10).Value = 40
$ws.Cells.Item(47, 11).Value = 2200
$ws.Cells.Item(47, 12).Value = 2200
$ws.Cells.Item(47, 13).Value = 2200
$ws.Cells.Item(47, 16).Value = 2200

# Row 48
$ws.Cells.Item(48, 4).Value = 44707
$ws.Cells.Item(48, 10).Value = 100
$ws.Cells.Item(48, 11).Value = 4700
$ws.Cells.Item(48, 12).Value = 4700
$ws.Cells.Item(48, 13).Value = 4700
$ws.Cells.Item(48, 16).Value = 4700

# Row 49
$ws.Cells.Item(49, 4).Value = 44784
$ws.Cells.Item(49, 10).Value = 180
$ws.Cells.Item(49, 11).Value = 2700
$ws.Cells.Item(49, 12).Value = 2700
$ws.Cells.Item(49, 13).Value = 2700
$ws.Cells.Item(49, 16).Value = 2700

# Row 50
$ws.Cells.Item(50, 4).Value = 44818
$ws.Cells.Item(50, 10).Value = 35
$ws.Cells.Item(50, 11).Value = 2700
$ws.Cells.Item(50, 12).Value = 2700
$ws.Cells.Item(50, 13).Value = 2700
$ws.Cells.Item(50, 16).Value = 2700

# Row 51
$ws.Cells.Item(51, 4).Value = 44881
$ws.Cells.Item(51, 10).Value = 100
$ws.Cells.Item(51, 11).Value = 2500
$ws.Cells.Item(51, 12).Value = 2500
$ws.Cells.Item(51, 13).Value = 2500
$ws.Cells.Item(51, 16).Value = 2500

# Row 52
$ws.Cells.Item(52, 4).Value = 44476
$ws.Cells.Item(52, 10).Value = 30
$ws.Cells.Item(52, 11).Value = 2200
$ws.Cells.Item(52, 12).Value = 2200
$ws.Cells.Item(52, 13).Value = 2200
$ws.Cells.Item(52, 16).Value = 2200

# Row 53
$ws.Cells.Item(53, 4).Value = 44798
$ws.Cells.Item(53, 10).Value = 80
$ws.Cells.Item(53, 11).Value = 2700
$ws.Cells.Item(53, 12).Value = 2700
$ws.Cells.Item(53, 13).Value = 2700
$ws.Cells.Item(53, 16).Value = 2700

# Row 54
$ws.Cells.Item(54, 4).Value = 44882
$ws.Cells.Item(54, 10).Value = 80

# Row 55
$ws.Cells.Item(55, 4).Value = 44832
$ws.Cells.Item(55, 10).Value = 80
$ws.Cells.Item(55, 11).Value = 2500
$ws.Cells.Item(55, 12).Value = 2500
$ws.Cells.Item(55, 13).Value = 2500
$ws.Cells.Item(55, 16).Value = 2500

# Row 56
$ws.Cells.Item(56, 4).Value = 44754
$ws.Cells.Item(56, 10).Value = 50
$ws.Cells.Item(56, 11).Value = 3300
$ws.Cells.Item(56, 12).Value = 3300
$ws.Cells.Item(56, 13).Value = 3300
$ws.Cells.Item(56, 16).Value = 3300

# Row 57
$ws.Cells.Item(57, 4).Value = 44855
$ws.Cells.Item(57, 10).Value = 30
$ws.Cells.Item(57, 11).Value = 2500
$ws.Cells.Item(57, 12).Value = 2500
$ws.Cells.Item(57, 13).Value = 2500
$ws.Cells.Item(57, 16).Value = 2500

# Row 58
$ws.Cells.Item(58, 4).Value = 44757
$ws.Cells.Item(58, 10).Value = 80
$ws.Cells.Item(58, 11).Value = 3300
$ws.Cells.Item(58, 12).Value = 3300
$ws.Cells.Item(58, 13).Value = 3300
$ws.Cells.Item(58, 16).Value = 3300

# Row 59
$ws.Cells.Item(59, 4).Value = 44795
$ws.Cells.Item(59, 10).Value = 120
$ws.Cells.Item(59, 11).Value = 2700
$ws.Cells.Item(59, 12).Value = 2700
$ws.Cells.Item(59, 13).Value = 2700
$ws.Cells.Item(59, 16).Value = 2700
